$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly record for "Espinaca" / Vega Monumental Concepción
# right before the current row 124. Excel shifts rows 124-128 down to
# 125-129 (and copies row 124's formatting, including the date style on
# column D) the same way a manual "Insert Row" above row 124 would.
$ws.Rows.Item(124).Insert()

$ws.Cells.Item(124, 1).Value = 11
$ws.Cells.Item(124, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(124, 3).Value = "Bíobío"
$ws.Cells.Item(124, 4).Value = 45106
$ws.Cells.Item(124, 5).Value = 8
$ws.Cells.Item(124, 6).Value = 100112012
$ws.Cells.Item(124, 7).Value = "Espinaca"
$ws.Cells.Item(124, 8).Value = "Sin especificar"
$ws.Cells.Item(124, 9).Value = "Primera"
$ws.Cells.Item(124, 10).Value = 50
$ws.Cells.Item(124, 11).Value = 7000
$ws.Cells.Item(124, 12).Value = 7500
$ws.Cells.Item(124, 13).Value = 7200
$ws.Cells.Item(124, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(124, 15).Value = "Región Metropolitana"
$ws.Cells.Item(124, 16).Value = 720
$ws.Cells.Item(124, 17).Value = 10
$ws.Cells.Item(124, 18).Value = "Hortaliza"
